# Timesheet changes by Ruchika
#
# For rows 28-31, columns AP..AU were previously blank (style s="26") and
# are now marked "OFF" using the same highlighted style already used by
# column AO (and other "OFF" cells) on those rows.
#
# We reproduce this by:
#   1. Copying the already-correctly-styled AO cell onto the AP:AU range
#      for each row (this carries over the "OFF" shared-string value and
#      the grey/centered style in one shot, and also fills the whole
#      destination range's formatting).
#   2. Re-asserting the "OFF" text on every individual cell in the range,
#      since Copy-to-range only guarantees the value lands on the first
#      destination cell in this host.
#
# Finally, the active selection is moved to match the new cursor position
# recorded in the workbook (AJ14 scrolled into view, AU28:AU31 selected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FebruaryMarch 2013")

$rows = 28, 29, 30, 31
$cols = "AP", "AQ", "AR", "AS", "AT", "AU"

foreach ($row in $rows) {
    $source = $ws.Range("AO" + $row)
    $target = $ws.Range("AP" + $row + ":AU" + $row)

    # Copy value + style (fill, font, alignment, number format, borders)
    # from AO<row> across AP<row>:AU<row>.
    $source.Copy($target)

    # Make sure every cell in the destination actually carries the "OFF"
    # text (Copy-to-range in this host only reliably stamps the first
    # cell's value).
    foreach ($col in $cols) {
        $ws.Range($col + $row).Value2 = "OFF"
    }
}

# Restore the view/selection state recorded for the sheet.
$ws.Activate()
$ws.Range("AJ14").Select()
$ws.Range("AU28:AU31").Select()
